$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to Text before writing, to avoid Excel
# auto-converting numeric-looking strings (e.g. "426.27") into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.646.08'
$ws.Range("E2").Value = '  +5.57%  '
$ws.Range("D3").Value = '3.850.58'
$ws.Range("E3").Value = '  +11.13%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '426.27'
$ws.Range("E5").Value = '  +9.89%  '
$ws.Range("D6").Value = '130.74'
$ws.Range("E6").Value = '  +6.94%  '
$ws.Range("D7").Value = '3.845.95'
$ws.Range("E7").Value = '  +7.58%  '
$ws.Range("D8").Value = '0.612'
$ws.Range("E8").Value = '  +4.97%  '
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").Value = '0.728'
$ws.Range("E10").Value = '  +9.23%  '
$ws.Range("D11").Value = '0.157'
$ws.Range("E11").Value = '  +10.15%  '
$ws.Range("D12").Value = '0.0000341'
$ws.Range("E12").Value = '  +8.23%  '
$ws.Range("D13").Value = '41.23'
$ws.Range("E13").Value = '  +7.25%  '
$ws.Range("D14").Value = '10.33'
$ws.Range("E14").Value = '  +13.27%  '
$ws.Range("D15").Value = '4.470.39'
$ws.Range("E15").Value = '  +13.10%  '
$ws.Range("D16").Value = '15.98'
$ws.Range("E16").Value = '  +27.95%  '
$ws.Range("D17").Value = '3.897.85'
$ws.Range("E17").Value = '  +12.86%  '
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("D19").Value = '19.98'
$ws.Range("E19").Value = '  +8.42%  '
$ws.Range("B20").Value = 'Polygon'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D20").Value = '1.09'
$ws.Range("E20").Value = '  +8.26%  '
$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").Value = '66.918.98'
$ws.Range("E21").Value = '  +6.33%  '
$ws.Range("D22").Value = '413.72'
$ws.Range("E22").Value = '  +6.23%  '
$ws.Range("D23").Value = '15.02'
$ws.Range("E23").Value = '  +9.35%  '
$ws.Range("D24").Value = '84.49'
$ws.Range("E24").Value = '  +6.02%  '
$ws.Range("D25").Value = '3.06'
$ws.Range("E25").Value = '  +8.12%  '
$ws.Range("D26").Value = '37.64'
$ws.Range("E26").Value = '  +14.99%  '
$ws.Range("D27").Value = '9.99'
$ws.Range("E27").Value = '  +13.80%  '
$ws.Range("E28").Value = '  +10.93%  '
$ws.Range("E29").Value = '  +2.90%  '
$ws.Range("D30").Value = '9.08'
$ws.Range("E30").Value = '  +35.54%  '
$ws.Range("D31").Value = '719.27'
$ws.Range("E31").Value = '  +8.79%  '
$ws.Range("D32").Value = '13.63'
$ws.Range("E32").Value = '  +16.01%  '
$ws.Range("E33").Value = '  +14.06%  '
$ws.Range("E34").Value = '  +6.93%  '
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("D36").Value = '39.19'
$ws.Range("E36").Value = '  +7.20%  '
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").Value = '55.78'
$ws.Range("E38").Value = '  +4.10%  '
$ws.Range("D39").Value = '5.47'
$ws.Range("E39").Value = '  +37.33%  '
$ws.Range("D40").Value = '0.0₃0748'
$ws.Range("E40").Value = '  +19.67%  '
$ws.Range("E41").Value = '  +6.67%  '
$ws.Range("D42").Value = '2.88'
$ws.Range("E42").Value = '  +9.17%  '
$ws.Range("E43").Value = '  +0.72%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").Value = '3.25'
$ws.Range("E44").Value = '  +8.34%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = '0.136'
$ws.Range("E45").Value = '  +4.13%  '
$ws.Range("B46").Value = 'LidoDAOToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D46").Value = '3.38'
$ws.Range("E46").Value = '  +11.10%  '
$ws.Range("D47").Value = '0.317'
$ws.Range("E47").Value = '  +15.65%  '
$ws.Range("D48").Value = '142.40'
$ws.Range("E48").Value = '  +2.92%  '
$ws.Range("D49").Value = '2.04'
$ws.Range("E49").Value = '  +5.61%  '
$ws.Range("D50").Value = '2.60'
$ws.Range("E50").Value = '  +6.16%  '
$ws.Range("E51").Value = '  +7.08%  '

# Restore default (unstyled) cell style for column D so the saved file
# does not pick up a spurious text-format style index.
$ws.Range("D2:D51").Style = "Normal"

